$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "theta_threshold_range" parameter row (row 5) is being removed entirely;
# deleting it shifts the "pie_threshold_range" row (old row 6) up into row 5.
$ws.Rows.Item(5).Delete()

# Row 6's "B" cell used a distinct (Times New Roman) style; match it back to
# the plain data style shared by the rest of the column by copying C5's
# formatting (already the right style) onto B5.
$ws.Range("C5").Copy()
$ws.Range("B5").PasteSpecial(-4122)

# Update the Min/Max threshold values per the new data.
$ws.Range("B2").Value = 3.8
$ws.Range("C2").Value = 12
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 9.6
$ws.Range("B4").Value = 0.8
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20

# Match the saved selection state seen in the target workbook.
$ws.Range("C3").Select()
